$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    # Long, already-unique strings: no need to restrict to whole words.
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

function Replace-WholeWord($old, $new) {
    # Short numeric tokens: restrict to whole-word matches. Find/Replace
    # with wdReplaceAll scans the *entire* document, so plain substring
    # search would also hit these numbers where they occur inside other,
    # longer numbers elsewhere in the invoice table (e.g. "92.0 " is a
    # substring of "492.0 "). Word boundaries also fall at table-cell
    # edges, so this still isolates each standalone cell correctly.
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# Invoice NO.
Replace-Unique "74231_2021731" "53672_202181"

# Date
Replace-Unique "2021-07-31" "2021-08-01"

# Line-item description
Replace-Unique "System wizyjny" "fsdfsd"

# Line-item QTY (4 -> 132)
Replace-WholeWord "4" "132"

# Line-item UNIT PRICE (100 -> 231)
Replace-WholeWord "100" "231"

# Tax Total (92.0 is also a substring of "492.0 ", so use whole-word)
Replace-WholeWord "92.0 " "37505.159999999996 "

# Balance Due (492.0 is also a substring of the new TOTAL value below)
Replace-WholeWord "492.0 " "67997.16 "

# Line-item TOTAL
Replace-WholeWord "400.0 " "30492.0 "

# Tax Rate (23 -> 123)
Replace-WholeWord "23" "123"
